# Update basic stats (allmus_size_vs_region_name) with new 2022 values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - small
$ws.Range("C2").Value = 70
$ws.Range("D2").Value = 420
$ws.Range("E2").Value = 164
$ws.Range("F2").Value = 165
$ws.Range("G2").Value = 274
$ws.Range("H2").Value = 166
$ws.Range("I2").Value = 52
$ws.Range("J2").Value = 156
$ws.Range("K2").Value = 335
$ws.Range("L2").Value = 333
$ws.Range("M2").Value = 146
$ws.Range("N2").Value = 133

# Row 3 - medium
$ws.Range("C3").Value = 20
$ws.Range("D3").Value = 118
$ws.Range("E3").Value = 62
$ws.Range("G3").Value = 91
$ws.Range("H3").Value = 69
$ws.Range("I3").Value = 39
$ws.Range("L3").Value = 136
$ws.Range("M3").Value = 84
$ws.Range("N3").Value = 85

# Row 4 - large
$ws.Range("C4").Value = 7
$ws.Range("D4").Value = 56
$ws.Range("E4").Value = 28
$ws.Range("F4").Value = 32
$ws.Range("G4").Value = 32
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 32
$ws.Range("J4").Value = 54
$ws.Range("N4").Value = 53

# Row 5 - huge
$ws.Range("H5").Value = 10

# Row 6 - unknown_sz
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 26
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 11
$ws.Range("G6").Value = 17
$ws.Range("H6").Value = 22
$ws.Range("I6").Value = 13
$ws.Range("J6").Value = 14
$ws.Range("K6").Value = 25
$ws.Range("L6").Value = 29
$ws.Range("M6").Value = 8
$ws.Range("N6").Value = 16
